$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add Link column, copy header style from B1 to C1
$ws.Range('C1').Value = 'Link'
$ws.Range('B1').Copy() | Out-Null
$ws.Range('C1').PasteSpecial(-4122) | Out-Null

# Row data: Name (A), Value (B), Link (C) for rows 2..51
$rows = @{}
$rows[2] = @('Base Carregador Controle Xbox Series S X Branco 2 Baterias', '124', 'https://click1.mercadolivre.com.br/mclics/clicks/external/MLB/count?a=IBWzvroGZTpsmo7HS6ecnn%2BE2DJ7H2Lh2vf6uoZuBZ%2FiDMdGFj5Zr87d3tEtjjTK7iSwv64iLjNmAunPU8v3uHB7vPE8DWRhbFlVEYXZYFKw%2B%2FpkgZsdx%2BaXF5RoU6nqAmqh%2Fas40l43qyETH37Uev8G%2FVkqMoqytJaF0mnT5QL789cAdl6dEjCa7IN8Cl2hjD6nbeToHlds7dRBocfM3y9vjAW3%2BcS4789FRFsqXzhZtNDdbbttm1%2FqOc%2BMspIDsNBDjqkE1aFTs4bOsZztANyfxL926tOnUFNCBlzVrFeb5cVgneIfubF62kZpxUJsi%2FPjzal3vhK5q6Nr%2BfRGrCjoXVKow%2BFKasivD56TUu8ibH%2FunwoqLILqgzWghd%2F%2FILHZVZCwwUneotaLetXf%2FCg1aE12%2FTvq0FxajSWp0wnqcxpz0pHr5qrFPoMWPagxGCFYCTzcoMdXJMe%2B7WzdUCzyHGQlUQoPEaDcuyRLgw%2Fol8OdidFUzVCLddVVsYd57Jd4XGHiRwjh6AtfyNXMJSyDi%2F1xOeTlcYOMdmnGnNlHMSR5hcyg4wGe7vdwMzxK2kcvAl%2FPn%2Fd7jqfYEDJh%2FvoJqkLUgmd%2FBaHcXiwRZwgwJyIe3TkbORrpVtWfEk1uIBYhsqr3w9SvnG4f0Yr8n0f5RJTO28KzMAEukFJ45zKxQccLOAqRwmbQDaaEgu%2FBCp5Px6U5V48sUmakhppgLHdeWRUdgAiDG1zvvP7zE%2FuabbhiJPx45y9VOTgdI0S%2FEhurDD2nk4C975mV4Q%3D%3D&rb=x#polycard_client=search-nordic&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54')
$rows[3] = @('Base Cooler Xbox Series S Carregador Controle Bateria 800mah', '198', 'https://click1.mercadolivre.com.br/mclics/clicks/external/MLB/count?a=coDdJ7DEks3kW6PblmiJLKya%2F8z%2FDoKC9MNh0w92xnG7X7RMeCVYKDp3tUI3En6NDwqeiKU%2Fi9EWj0evKhC4UUxVKNoTuGkextgianGy5LFe0Jj61o279LX23Hi0wbTyl8DA2PMzKjZXe4RhoMGp9eTigDr6smy9bUHjwAyXBo9VwyY1R9m6z22FeS69fVOI4ND%2BNHPdO%2F1WOPWuBdpoFVaYUfDAupSLGbR7jHRq5X4vUNBGF4newNQBMHq9wPLzax4slmuFL2%2BBkU72zGg16ZksRo%2Bdl23SUoilJOwiYyPfmuWbj1iSAdbfmJfzZS966IUDv6GPHyrFGqivLdX%2BWM%2FhjFzoHGU62eHVQ6djcKbOoEsdEaBDQ0FoeUkXIR69k1g8Bjc5JCWmNAYKrSSW14f99udTkKfB62FizR6xWACO917eeWgho%2Bcfr3ZIJkPYZ03fj7ldGGuS2%2BEv7nukqc66PRyqD7bnT0bzJT7HJ5ZreBoGeOsCIT4Q62DPtbGdpM8SIhdd1Z4nwcK%2Fhh6zOj%2BqOfnyZa860TnlHpFfCWfIB%2B3NfP9hp4L1Cw5Md4I8H9Hccnw70DZmuap9KvPgzgRE8K6xpVVh2a5EqaC5VuAs1TW5acUj0h%2FHSEvPsgFCPhaI9zOrLhIOmdkVJWS4JQfmjfRB%2FWlENE7n6ooy6ngj5EMsMWLA4Va%2BID7uZVecGG9ICd4OFrWI83fxPd%2F%2FWhJ3lF4QOYErTHSCuqmYVHG8g1c1Bt1LhofjYyEF%2FJLtyaVpGvbhcH8mot0W9y5%2FqbY%3D&rb=x#polycard_client=search-nordic&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54')
$rows[4] = @('Controle Xbox Wireless Series X|s Carbon black', '629', 'https://www.mercadolivre.com.br/controle-xbox-wireless-series-xs-carbon-black/p/MLB16268160#polycard_client=search-nordic&searchVariation=MLB16268160&wid=MLB3998223345&position=4&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[5] = @('Controle Microsoft Xbox Wireless Series X/s Velocity Green Cor Verde', '545', 'https://www.mercadolivre.com.br/controle-microsoft-xbox-wireless-series-xs-velocity-green-cor-verde/p/MLB22538792#polycard_client=search-nordic&searchVariation=MLB22538792&wid=MLB3443874399&position=3&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[6] = @('Controle joystick sem fio Microsoft Xbox Elite Series 2 branco', '1.080', 'https://www.mercadolivre.com.br/controle-joystick-sem-fio-microsoft-xbox-elite-series-2-branco/p/MLB19627127#polycard_client=search-nordic&searchVariation=MLB19627127&wid=MLB3751953327&position=6&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[7] = @('Controle joystick sem fio Microsoft Wireless Controller Series X|S Series X e S electric volt', '430', 'https://www.mercadolivre.com.br/controle-joystick-sem-fio-microsoft-wireless-controller-series-xs-series-x-e-s-electric-volt/p/MLB18030215#polycard_client=search-nordic&searchVariation=MLB18030215&wid=MLB2005787418&position=7&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[8] = @('Base Carregador Controle Xbox Series S X Branco 2 Baterias', '124', 'https://produto.mercadolivre.com.br/MLB-3361773385-base-carregador-controle-xbox-series-s-x-branco-2-baterias-_JM#polycard_client=search-nordic&position=32&search_layout=grid&type=item&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54')
$rows[9] = @('Controle joystick sem fio Microsoft Xbox Carbon Black preto', '409', 'https://www.mercadolivre.com.br/controle-joystick-sem-fio-microsoft-xbox-carbon-black-preto/p/MLB17483958#polycard_client=search-nordic&searchVariation=MLB17483958&wid=MLB3181675322&position=5&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[10] = @('Controle joystick sem fio Microsoft Wireless Controller Series X|S Series X e S deep pink', '499', 'https://www.mercadolivre.com.br/controle-joystick-sem-fio-microsoft-wireless-controller-series-xs-series-x-e-s-deep-pink/p/MLB22522892#polycard_client=search-nordic&searchVariation=MLB22522892&wid=MLB3584569802&position=13&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[11] = @('Controle Joystick Sem Fio Microsoft Xbox X/s Branco', '428', 'https://www.mercadolivre.com.br/controle-joystick-sem-fio-microsoft-xbox-xs-branco/p/MLB22225047#polycard_client=search-nordic&searchVariation=MLB22225047&wid=MLB5273289302&position=16&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[12] = @('Knup GM019 Controle Com Fio Usb Para Computador Rgb Gamepass', '129', 'https://www.mercadolivre.com.br/knup-gm019-controle-com-fio-usb-para-computador-rgb-gamepass/p/MLB21723314#polycard_client=search-nordic&searchVariation=MLB21723314&wid=MLB4454723210&position=10&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[13] = @('Para Peças Acessórias De Controle Elite Series 2, Jogo 13 Em', '163', 'https://www.mercadolivre.com.br/for-elite-series-2-controller-accessory-parts-13-in-1-game/p/MLB2001851161#polycard_client=search-nordic&searchVariation=MLB2001851161&wid=MLB3866470286&position=9&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[14] = @('Base Cooler Xbox Series S Carregador Controle Bateria 800mah', '198', 'https://produto.mercadolivre.com.br/MLB-3259868703-base-cooler-xbox-series-s-carregador-controle-bateria-800mah-_JM#polycard_client=search-nordic&position=33&search_layout=grid&type=item&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54')
$rows[15] = @('Controle Sem Fio Xbox Wireless Preto Cor Carbon black', '414', 'https://www.mercadolivre.com.br/controle-sem-fio-xbox-wireless-preto-cor-carbon-black/p/MLB23280037#polycard_client=search-nordic&searchVariation=MLB23280037&wid=MLB5181374502&position=17&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[16] = @('Controle Sem Fio Microsoft Xbox Wireless Series X|s Cor Carbon black', '410', 'https://www.mercadolivre.com.br/controle-sem-fio-microsoft-xbox-wireless-series-xs-cor-carbon-black/p/MLB26089834#polycard_client=search-nordic&searchVariation=MLB26089834&wid=MLB3997678171&position=11&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[17] = @('Kontrol Freek Extensor Analógico Grip Control Fps Freek', '65', 'https://produto.mercadolivre.com.br/MLB-4977470558-kontrol-freek-extensor-analogico-grip-control-fps-freek-_JM?searchVariation=181291740168#polycard_client=search-nordic&searchVariation=181291740168&position=34&search_layout=grid&type=item&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54')
$rows[18] = @('Controle Sem Fio Xbox One E Cabo Usb Preto', '479', 'https://produto.mercadolivre.com.br/MLB-4958466016-controle-sem-fio-xbox-one-e-cabo-usb-preto-_JM?searchVariation=181222590974#polycard_client=search-nordic&searchVariation=181222590974&position=35&search_layout=grid&type=item&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54')
$rows[19] = @('Controle joystick sem fio Microsoft Xbox Wireless Controller Series X|S Series X e S robot white', '414', 'https://www.mercadolivre.com.br/controle-joystick-sem-fio-microsoft-xbox-wireless-controller-series-xs-series-x-e-s-robot-white/p/MLB16268161#polycard_client=search-nordic&searchVariation=MLB16268161&wid=MLB3492668122&position=15&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[20] = @('2 Botões Gatilho Lb Rb Para Controle Xbox Séries X S 1914', '29', 'https://produto.mercadolivre.com.br/MLB-4831376386-2-botoes-gatilho-lb-rb-para-controle-xbox-series-x-s-1914-_JM?searchVariation=183399503185#polycard_client=search-nordic&searchVariation=183399503185&position=36&search_layout=grid&type=item&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54')
$rows[21] = @('Controle Microsoft Xbox One Series X E S', '419', 'https://www.mercadolivre.com.br/controle-microsoft-xbox-one-series-x-e-s/p/MLB23097417#polycard_client=search-nordic&searchVariation=MLB23097417&wid=MLB4010266483&position=20&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[22] = @('Botão Gatilho Rb Lb  Para Controle Xbox Series S E X', '16', 'https://produto.mercadolivre.com.br/MLB-3755920903-boto-gatilho-rb-lb-para-controle-xbox-series-s-e-x-_JM?searchVariation=183324883455#polycard_client=search-nordic&searchVariation=183324883455&position=37&search_layout=grid&type=item&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54')
$rows[23] = @('Joystick Microsoft XBox One v1 Preto sem fio', '449', 'https://www.mercadolivre.com.br/joystick-microsoft-xbox-one-v1-preto-sem-fio/p/MLB12384031#polycard_client=search-nordic&searchVariation=MLB12384031&wid=MLB4768575412&position=14&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[24] = @('Controle Xbox 360 Sem Fio', '129', 'https://www.mercadolivre.com.br/controle-xbox-360-sem-fio/p/MLB24648019#polycard_client=search-nordic&searchVariation=MLB24648019&wid=MLB3626412625&position=8&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[25] = @('Carregador Controle Xbox One Series S X + 2 Baterias 800mah', '139', 'https://produto.mercadolivre.com.br/MLB-2712532600-carregador-controle-xbox-one-series-s-x-2-baterias-800mah-_JM?searchVariation=178400028139#polycard_client=search-nordic&searchVariation=178400028139&position=38&search_layout=grid&type=item&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54')
$rows[26] = @('Controle joystick sem fio Microsoft Wireless Controller Series X|S Series X e S pulse red', '448', 'https://www.mercadolivre.com.br/controle-joystick-sem-fio-microsoft-wireless-controller-series-xs-series-x-e-s-pulse-red/p/MLB17375584#polycard_client=search-nordic&searchVariation=MLB17375584&wid=MLB5284364390&position=21&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[27] = @('Control Freek - Fps Freek Galaxy Xbox Série X/s - Um', '46', 'https://www.mercadolivre.com.br/control-freek-fps-freek-galaxy-xbox-serie-xs-um/p/MLB20883537#polycard_client=search-nordic&searchVariation=MLB20883537&wid=MLB3798096643&position=24&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[28] = @('Analógico Controle Xbox One Series S E X Direcional Original', '16', 'https://produto.mercadolivre.com.br/MLB-2794264431-analogico-controle-xbox-one-series-s-e-x-direcional-original-_JM#polycard_client=search-nordic&position=39&search_layout=grid&type=item&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54')
$rows[29] = @('Controle Compatível Xbox One Series E Pc C/ Fio Top', '139', 'https://www.mercadolivre.com.br/controle-compativel-xbox-one-series-e-pc-c-fio-top/p/MLB22366874#polycard_client=search-nordic&searchVariation=MLB22366874&wid=MLB3856360757&position=12&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[30] = @('1 Botão Gatilho Rb Lb  Para Controle Xbox Series S E X Preto', '16', 'https://produto.mercadolivre.com.br/MLB-3756445547-1-boto-gatilho-rb-lb-para-controle-xbox-series-s-e-x-preto-_JM?searchVariation=183339333763#polycard_client=search-nordic&searchVariation=183339333763&position=40&search_layout=grid&type=item&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54')
$rows[31] = @('Bateria Controle Xbox Séries S X 1200mah Cabo Carregador 3m', '59', 'https://produto.mercadolivre.com.br/MLB-5111737986-bateria-controle-xbox-series-s-x-1200mah-cabo-carregador-3m-_JM#polycard_client=search-nordic&position=41&search_layout=grid&type=item&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54')
$rows[32] = @('2 Botões Do Analógico Para Controle Xbox One Slim Series S', '14', 'https://produto.mercadolivre.com.br/MLB-4862989250-2-botoes-do-analogico-para-controle-xbox-one-slim-series-s-_JM?searchVariation=180963648648#polycard_client=search-nordic&searchVariation=180963648648&position=42&search_layout=grid&type=item&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54')
$rows[33] = @('Controle sem fio Xbox Elite Series 2 Core Red', '1.198', 'https://www.mercadolivre.com.br/controle-sem-fio-xbox-elite-series-2-core-red/p/MLB22892331#polycard_client=search-nordic&searchVariation=MLB22892331&wid=MLB3859552111&position=23&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[34] = @('Joystick sem fio Microsoft Xbox Y PC Gamepad branco', '469', 'https://www.mercadolivre.com.br/joystick-sem-fio-microsoft-xbox-y-pc-gamepad-branco/p/MLB37263008#polycard_client=search-nordic&searchVariation=MLB37263008&wid=MLB5281114478&position=19&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[35] = @('Analógico Controle Para Xbox One Series S E X Direcional ', '14', 'https://produto.mercadolivre.com.br/MLB-3445580683-analogico-controle-para-xbox-one-series-s-e-x-direcional-_JM#polycard_client=search-nordic&position=43&search_layout=grid&type=item&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54')
$rows[36] = @('Controle Joystick Sem Fio Microsoft Xbox Xbox Series ', '425', 'https://produto.mercadolivre.com.br/MLB-4474529196-controle-joystick-sem-fio-microsoft-xbox-xbox-series-_JM?searchVariation=179921638214#polycard_client=search-nordic&searchVariation=179921638214&position=44&search_layout=grid&type=item&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54')
$rows[37] = @('Controle Sem Fio Xbox Sky Cipher Special Edition', '528', 'https://www.mercadolivre.com.br/controle-sem-fio-xbox-sky-cipher-special-edition/p/MLB38804309#polycard_client=search-nordic&searchVariation=MLB38804309&wid=MLB5196273196&position=26&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[38] = @('Cabo Carregador Para Controle Xbox Series Tipo C 3 Metros', '36', 'https://produto.mercadolivre.com.br/MLB-2223135747-cabo-carregador-para-controle-xbox-series-tipo-c-3-metros-_JM#polycard_client=search-nordic&position=45&search_layout=grid&type=item&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54')
$rows[39] = @('Controle 8bitdo Xbox Ultimate C Com Hall Effect A P/ Cor Cinza Escuro', '336', 'https://www.mercadolivre.com.br/controle-8bitdo-xbox-ultimate-c-com-hall-effect-a-p-cor-cinza-escuro/p/MLB37637160#polycard_client=search-nordic&searchVariation=MLB37637160&wid=MLB4002623151&position=18&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[40] = @('Controle Joystick Sem Fio Microsoft Xbox Shock Blue', '519', 'https://www.mercadolivre.com.br/controle-joystick-sem-fio-microsoft-xbox-shock-blue/p/MLB16268159#polycard_client=search-nordic&searchVariation=MLB16268159&wid=MLB3610784199&position=30&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[41] = @('Capa Silicon Control Xbox Series X S para alças de borracha preta', '23', 'https://www.mercadolivre.com.br/capa-silicon-control-xbox-series-x-s-para-alcas-de-borracha-preta/p/MLB34261805#polycard_client=search-nordic&searchVariation=MLB34261805&wid=MLB4011719819&position=27&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[42] = @('Controle Joystick Sem Fio Microsoft Xbox Xbox Series X|s Con', '454', 'https://produto.mercadolivre.com.br/MLB-3601791199-controle-joystick-sem-fio-microsoft-xbox-xbox-series-xs-con-_JM?searchVariation=181782924599#polycard_client=search-nordic&searchVariation=181782924599&position=46&search_layout=grid&type=item&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54')
$rows[43] = @('Controle Sem Fio Microsoft Xbox Xt7 Pro Carbon Black', '430', 'https://produto.mercadolivre.com.br/MLB-3986451857-controle-sem-fio-microsoft-xbox-xt7-pro-carbon-black-_JM?searchVariation=182998549084#polycard_client=search-nordic&searchVariation=182998549084&position=47&search_layout=grid&type=item&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54')
$rows[44] = @('Controle Joystick Microsoft Xbox Series X/s White Robot', '415', 'https://www.mercadolivre.com.br/controle-joystick-microsoft-xbox-series-xs-white-robot/p/MLB27744165#polycard_client=search-nordic&searchVariation=MLB27744165&wid=MLB5309401862&position=31&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[45] = @('Controle Sem Fio Xbox Series S X One Pc Original Preto', '479', 'https://produto.mercadolivre.com.br/MLB-3397445599-controle-sem-fio-xbox-series-s-x-one-pc-original-preto-_JM?searchVariation=179421230071#polycard_client=search-nordic&searchVariation=179421230071&position=48&search_layout=grid&type=item&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54')
$rows[46] = @('2 Analógicos 3d Controle + T8 Para Xbox One / One S / Séries', '37', 'https://produto.mercadolivre.com.br/MLB-2774190133-2-analogicos-3d-controle-t8-para-xbox-one-one-s-series-_JM#polycard_client=search-nordic&position=49&search_layout=grid&type=item&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54')
$rows[47] = @('Controle Joystick Sem Fio Xbox Series X|s Con', '449', 'https://produto.mercadolivre.com.br/MLB-4474553990-controle-joystick-sem-fio-xbox-series-xs-con-_JM?searchVariation=181872070455#polycard_client=search-nordic&searchVariation=181872070455&position=50&search_layout=grid&type=item&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54')
$rows[48] = @('Controle Microsoft Starfield Edition Series X|S, One, Pc', '1.597', 'https://www.mercadolivre.com.br/controle-microsoft-starfield-edition-series-xs-one-pc/p/MLB24045273#polycard_client=search-nordic&searchVariation=MLB24045273&wid=MLB3763294656&position=22&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[49] = @('Bateria Recarregável Para Controle Xbox One Tampa Traseira C Cor Preto', '23', 'https://www.mercadolivre.com.br/bateria-recarregavel-para-controle-xbox-one-tampa-traseira-c-cor-preto/p/MLB25898760#polycard_client=search-nordic&searchVariation=MLB25898760&wid=MLB3897236370&position=29&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[50] = @('Controle joystick sem fio Microsoft Slim One S branco', '429', 'https://www.mercadolivre.com.br/controle-joystick-sem-fio-microsoft-slim-one-s-branco/p/MLB21685852#polycard_client=search-nordic&searchVariation=MLB21685852&wid=MLB3868783357&position=25&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')
$rows[51] = @('Controle Joystick Sem Fio Microsoft Xbox Wireless Series X|s Cor Shock Blue', '519', 'https://www.mercadolivre.com.br/controle-joystick-sem-fio-microsoft-xbox-wireless-series-xs-cor-shock-blue/p/MLB39972584#polycard_client=search-nordic&searchVariation=MLB39972584&wid=MLB5135365484&position=28&search_layout=grid&type=product&tracking_id=42fdeefd-8fa5-49e5-99ce-87042b654a54&sid=search')

foreach ($r in $rows.Keys) {
  $vals = $rows[$r]
  $ws.Cells.Item($r, 1).Value = $vals[0]
  # Column B holds numeric-looking text (Brazilian formatted prices); force text storage
  $bCell = $ws.Cells.Item($r, 2)
  $bCell.NumberFormat = '@'
  $bCell.Value = $vals[1]
  $bCell.ClearFormats()
  $ws.Cells.Item($r, 3).Value = $vals[2]
}

Write-Output 'done'